# Amended the meeting diary
#
# Edits the Sheet2 "meeting diary" log:
#  - Changes the first logged meeting's date
#  - Adds a brand new meeting entry (row 8) that replaces the old
#    placeholder "Fill in as needed" row, recording a catch-up meeting
#    with the previous team about the wild sighting data
#  - Updates the group members cell to add a space after the commas

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- First meeting row (row 7): date moves from 14 Sep to 14 Oct 2023 ---
$ws2.Range("A7").Value2 = 45213

# --- New meeting row (row 8), replacing the old "Fill in as needed" row ---
$ws2.Range("A8").Value2 = 45215
$ws2.Range("B8").Value2 = 0.83333333333333337
$ws2.Range("C8").Value2 = 0.85416666666666663
# Leading "-" needs a quote prefix so Excel treats it as text, not a formula
$ws2.Range("D8").Formula = "'- All`n- Previous Team"
$ws2.Range("E8").Value2 = "Getting insights from the previous team who handle the wild sighting"

# --- Group members cell (B2): add spacing after commas ---
$ws2.Range("B2").Value2 = "Evan, Amrita, Phuong Mai"

# Copy the date/time number formats down from row 7 into the new row 8 cells
$ws2.Range("A7").Copy() | Out-Null
$ws2.Range("A8").PasteSpecial(-4122) | Out-Null
$ws2.Range("B7:C7").Copy() | Out-Null
$ws2.Range("B8:C8").PasteSpecial(-4122) | Out-Null
$ws2.Application.CutCopyMode = $false

# Wrap the new discussion text, then carry that formatting across to E8
$ws2.Range("D8").WrapText = $true
$ws2.Range("D8").Copy() | Out-Null
$ws2.Range("E8").PasteSpecial(-4122) | Out-Null
$ws2.Application.CutCopyMode = $false

# Rows 7 and 8 resize to fit the (now two-line) wrapped text
$ws2.Range("A7").EntireRow.RowHeight = 31
$ws2.Range("A8").EntireRow.RowHeight = 31

# Restore the active selection on Sheet2
$ws2.Activate() | Out-Null
$ws2.Range("C3").Select() | Out-Null
